$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1201.2222
$ws.Range("I19").Value = 1395.2
$ws.Range("K19").Value = 1395.2
$ws.Range("M19").Value = -1220.2
$ws.Range("H32").Value = 2992
$ws.Range("I32").Value = 3241
$ws.Range("J32").Value = 2707.4285
$ws.Range("K32").Value = 3241
$ws.Range("L32").Value = 2707.4285
$ws.Range("M32").Value = -2915
$ws.Range("N32").Value = -3359.4285
$ws.Range("H132").Value = 2080.125
$ws.Range("I132").Value = 935.55554
$ws.Range("K132").Value = 2806.66662
$ws.Range("M132").Value = -276.66662

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H76").Value = 15000
$ws.Range("J76").Value = 15000
$ws.Range("L76").Value = 15000
$ws.Range("N76").Value = -15676
$ws.Range("H79").Value = 15000
$ws.Range("J79").Value = 15000
$ws.Range("L79").Value = 15000
$ws.Range("N79").Value = -17340
$ws.Range("H82").Value = 38000
$ws.Range("J82").Value = 38000
$ws.Range("L82").Value = 38000
$ws.Range("N82").Value = -38722
$ws.Range("H85").Value = 38000
$ws.Range("J85").Value = 38000
$ws.Range("L85").Value = 38000
$ws.Range("N85").Value = -40496

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1944.5333
$ws.Range("I99").Value = 2383.625
$ws.Range("J99").Value = 1442.7142
$ws.Range("K99").Value = 2383.625
$ws.Range("L99").Value = 1442.7142
$ws.Range("M99").Value = -885.625
$ws.Range("N99").Value = -4438.7142

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 480.2
$ws.Range("I22").Value = 300.5
$ws.Range("J22").Value = 600
$ws.Range("K22").Value = 300.5
$ws.Range("L22").Value = 600
$ws.Range("M22").Value = 49.5
$ws.Range("N22").Value = -1300
$ws.Range("H31").Value = 2103.3635
$ws.Range("I31").Value = 1786
$ws.Range("J31").Value = 3282.1428
$ws.Range("K31").Value = 1786
$ws.Range("L31").Value = 3282.1428
$ws.Range("M31").Value = -1491
$ws.Range("N31").Value = -3872.1428
$ws.Range("H34").Value = 2103.3635
$ws.Range("I34").Value = 1786
$ws.Range("J34").Value = 3282.1428
$ws.Range("K34").Value = 1786
$ws.Range("L34").Value = 3282.1428
$ws.Range("M34").Value = -1584
$ws.Range("N34").Value = -3686.1428
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H74").Value = 14025.6
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 14025.6
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 14025.6
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -15773.6
$ws.Range("H77").Value = 14025.6
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 14025.6
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 42076.8
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -50812.8
$ws.Range("H80").Value = 19200
$ws.Range("J80").Value = 19200
$ws.Range("L80").Value = 19200
$ws.Range("N80").Value = -21446
$ws.Range("H81").Value = 48333.332
$ws.Range("J81").Value = 48333.332
$ws.Range("L81").Value = 48333.332
$ws.Range("N81").Value = -50329.332
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("M82").ClearContents()
$ws.Range("N82").ClearContents()
$ws.Range("H83").Value = 19200
$ws.Range("J83").Value = 19200
$ws.Range("L83").Value = 57600
$ws.Range("N83").Value = -68832
$ws.Range("H84").Value = 48333.332
$ws.Range("J84").Value = 48333.332
$ws.Range("L84").Value = 144999.996
$ws.Range("N84").Value = -154983.996
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("M85").ClearContents()
$ws.Range("N85").ClearContents()
$ws.Range("H87").Value = 29900
$ws.Range("J87").Value = 29900
$ws.Range("L87").Value = 29900
$ws.Range("N87").Value = -32272
$ws.Range("H90").Value = 29900
$ws.Range("J90").Value = 29900
$ws.Range("L90").Value = 89700
$ws.Range("N90").Value = -101556
$ws.Range("H132").Value = 2689.2727
$ws.Range("I132").Value = 717
$ws.Range("J132").Value = 4332.8335
$ws.Range("K132").Value = 2151
$ws.Range("L132").Value = 12998.5005
$ws.Range("M132").Value = 379
$ws.Range("N132").Value = -18058.5005
$ws.Range("H134").Value = 1895.1177
$ws.Range("I134").Value = 1483.5834
$ws.Range("J134").Value = 2882.8
$ws.Range("K134").Value = 4450.7502
$ws.Range("L134").Value = 8648.400000000001
$ws.Range("M134").Value = -1915.7502
$ws.Range("N134").Value = -13718.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 41031.64
$ws.Range("I4").Value = 217.68182
$ws.Range("K4").Value = 653.0454599999999
$ws.Range("M4").Value = -541.0454599999999
$ws.Range("H35").Value = 1499.5
$ws.Range("J35").Value = 1499.5
$ws.Range("L35").Value = 4498.5
$ws.Range("N35").Value = -5074.5
$ws.Range("H37").Value = 87662.336
$ws.Range("J37").Value = 87662.336
$ws.Range("L37").Value = 262987.008
$ws.Range("N37").Value = -263211.008
$ws.Range("H38").Value = 105.181816
$ws.Range("I38").Value = 67.5
$ws.Range("J38").Value = 150.4
$ws.Range("K38").Value = 202.5
$ws.Range("L38").Value = 451.2
$ws.Range("M38").Value = 144.5
$ws.Range("N38").Value = -1145.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1209.5333
$ws.Range("I16").Value = 1209.5333
$ws.Range("K16").Value = 1209.5333
$ws.Range("M16").Value = -1039.5333
$ws.Range("H22").Value = 523.0909
$ws.Range("I22").Value = 624.25
$ws.Range("J22").Value = 401.7
$ws.Range("K22").Value = 624.25
$ws.Range("L22").Value = 401.7
$ws.Range("M22").Value = -329.25
$ws.Range("N22").Value = -991.7
$ws.Range("H27").Value = 523.0909
$ws.Range("I27").Value = 624.25
$ws.Range("J27").Value = 401.7
$ws.Range("K27").Value = 624.25
$ws.Range("L27").Value = 401.7
$ws.Range("M27").Value = -517.25
$ws.Range("N27").Value = -615.7
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
$ws.Range("H40").Value = 2446.4614
$ws.Range("I40").Value = 2477.6667
$ws.Range("J40").Value = 2376.25
$ws.Range("K40").Value = 2477.6667
$ws.Range("L40").Value = 2376.25
$ws.Range("M40").Value = -2341.6667
$ws.Range("N40").Value = -2648.25
$ws.Range("H46").Value = 739.34485
$ws.Range("I46").Value = 836.05
$ws.Range("J46").Value = 524.44446
$ws.Range("K46").Value = 836.05
$ws.Range("L46").Value = 524.44446
$ws.Range("M46").Value = -648.05
$ws.Range("N46").Value = -900.44446
$ws.Range("H55").Value = 148.36734
$ws.Range("I55").Value = 132.58974
$ws.Range("J55").Value = 209.9
$ws.Range("K55").Value = 132.58974
$ws.Range("L55").Value = 209.9
$ws.Range("M55").Value = 40.41025999999999
$ws.Range("N55").Value = -555.9
$ws.Range("H122").Value = 2683.3572
$ws.Range("I122").Value = 2389.125
$ws.Range("J122").Value = 4448.75
$ws.Range("K122").Value = 7167.375
$ws.Range("L122").Value = 13346.25
$ws.Range("M122").Value = -4717.375
$ws.Range("N122").Value = -18246.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3768.16
$ws.Range("I122").Value = 3529.95
$ws.Range("K122").Value = 10589.85
$ws.Range("M122").Value = -8139.849999999999
$ws.Range("H136").Value = 1547.85
$ws.Range("I136").Value = 1188.6154
$ws.Range("J136").Value = 2215
$ws.Range("K136").Value = 3565.8462
$ws.Range("L136").Value = 6645
$ws.Range("M136").Value = -1015.8462
$ws.Range("N136").Value = -11745
